# Charlotte Hornets 2023 roster - reorder rows 9-16 (Mark Williams, LaMelo Ball,
# Theo Maledon, Kai Jones, Bryce McGowens, James Bouknight, Cody Martin,
# Svi Mykhailiuk) into their new order (Kai Jones, Mark Williams, Bryce
# McGowens, LaMelo Ball, Theo Maledon, James Bouknight, Svi Mykhailiuk,
# Cody Martin), per the commit's restructuring.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 9-16, columns C..K (A and B... B changes too, A stays 7..14)
# Columns: B=No. C=Player D=Pos E=Ht F=Wt G=BirthDate H=Nat(country code) I=Exp J=College K=URL

$rowsData = @(
    @{ Row=9;  B=23; C="Kai Jones";          D="C";  E="6-11"; F=218; G="January 19, 2001";   H="bs"; I="1"; J="Texas";            K="https://www.basketball-reference.com/players/j/joneska01.html" },
    @{ Row=10; B=5;  C="Mark Williams";      D="C";  E="7-1";  F=241; G="December 16, 2001";  H="us"; I="R"; J="Duke";             K="https://www.basketball-reference.com/players/w/willima07.html" },
    @{ Row=11; B=7;  C="Bryce McGowens";     D="SG"; E="6-7";  F=179; G="November 8, 2002";   H="us"; I="R"; J="Nebraska";         K="https://www.basketball-reference.com/players/m/mcgowbr01.html" },
    @{ Row=12; B=1;  C="LaMelo Ball";        D="PG"; E="6-7";  F=180; G="August 22, 2001";    H="us"; I="2"; J="";                 K="https://www.basketball-reference.com/players/b/ballla01.html" },
    @{ Row=13; B=9;  C="Théo Maledon (TW)";  D="PG"; E="6-4";  F=175; G="June 12, 2001";       H="fr"; I="2"; J="";                 K="https://www.basketball-reference.com/players/m/maledth01.html" },
    @{ Row=14; B=2;  C="James Bouknight";    D="SG"; E="6-5";  F=190; G="September 18, 2000"; H="us"; I="1"; J="UConn";            K="https://www.basketball-reference.com/players/b/bouknja01.html" },
    @{ Row=15; B=10; C="Svi Mykhailiuk";     D="SF"; E="6-7";  F=205; G="June 10, 1997";       H="ua"; I="4"; J="Kansas";           K="https://www.basketball-reference.com/players/m/mykhasv01.html" },
    @{ Row=16; B=11; C="Cody Martin";        D="SF"; E="6-5";  F=205; G="September 28, 1995"; H="us"; I="3"; J="NC State, Nevada"; K="https://www.basketball-reference.com/players/m/martico01.html" }
)

foreach ($entry in $rowsData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F
    $ws.Cells.Item($r, 7).Value = $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
    $ws.Cells.Item($r, 11).Value = $entry.K
}

# Hyperlinks: in this COM engine, deleting hyperlinks on any single range clears
# ALL hyperlinks on the sheet, so rebuild the full K2:K16 hyperlink set in the
# correct (new) row order in one pass.
$ws.Range("K9").Hyperlinks.Delete()

$allUrls = @(
    @{ Row=2;  Url="https://www.basketball-reference.com/players/w/washipj01.html" },
    @{ Row=3;  Url="https://www.basketball-reference.com/players/r/roziete01.html" },
    @{ Row=4;  Url="https://www.basketball-reference.com/players/t/thorjt01.html" },
    @{ Row=5;  Url="https://www.basketball-reference.com/players/r/richani01.html" },
    @{ Row=6;  Url="https://www.basketball-reference.com/players/s/smithde03.html" },
    @{ Row=7;  Url="https://www.basketball-reference.com/players/o/oubreke01.html" },
    @{ Row=8;  Url="https://www.basketball-reference.com/players/h/haywago01.html" },
    @{ Row=9;  Url="https://www.basketball-reference.com/players/j/joneska01.html" },
    @{ Row=10; Url="https://www.basketball-reference.com/players/w/willima07.html" },
    @{ Row=11; Url="https://www.basketball-reference.com/players/m/mcgowbr01.html" },
    @{ Row=12; Url="https://www.basketball-reference.com/players/b/ballla01.html" },
    @{ Row=13; Url="https://www.basketball-reference.com/players/m/maledth01.html" },
    @{ Row=14; Url="https://www.basketball-reference.com/players/b/bouknja01.html" },
    @{ Row=15; Url="https://www.basketball-reference.com/players/m/mykhasv01.html" },
    @{ Row=16; Url="https://www.basketball-reference.com/players/m/martico01.html" }
)

foreach ($entry in $allUrls) {
    $cell = $ws.Cells.Item($entry.Row, 11)
    $ws.Hyperlinks.Add($cell, $entry.Url)
}
